$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new order line (row 12) for the Sandwich Picks glob-order lookup fix.
# Columns C/D/E hold numeric-looking values ("2", "11.81", "23.62") that must
# be stored as text (matching every other data cell in this sheet), so force
# a text number format on those three cells before assigning their values -
# otherwise Excel auto-coerces them into real numbers.
$ws.Range("C12:E12").NumberFormat = "@"

$ws.Range("A12").Value = "R801"
$ws.Range("B12").Value = "Sandwich Picks"
$ws.Range("C12").Value = "2"
$ws.Range("D12").Value = "11.81"
$ws.Range("E12").Value = "23.62"
